# Apply "Add data for 2022-12-01" update to carjacking-by-month-yoy-historical.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (Through 2022-11-22 -> Through 2022-11-23)
$ws.Name = "Through 2022-11-23"

# Update the header label for the "through" column (column I, row 1)
$ws.Range("I1").Value = "2022 (through 11-23)"

# Update the year-to-date totals for October and November, and the grand total
$ws.Range("I11").Value = 124   # October
$ws.Range("I12").Value = 87    # November
$ws.Range("I14").Value = 1485  # Grand total
